$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "40.308.33"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.221.94"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "297.51"
$ws.Range("E5").Value = "  +2.17%  "

$ws.Range("D6").Value = "88.37"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("D7").Value = "0.515"
$ws.Range("E7").Value = "  +0.58%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("D10").Value = "52.84"
$ws.Range("E10").Value = "  +8.42%  "

$ws.Range("D11").Value = "31.15"
$ws.Range("E11").Value = "  +2.35%  "

$ws.Range("D12").Value = "0.0784"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "2.557.61"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "13.93"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").Value = "2.199.33"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "0.740"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "40.239.06"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "11.41"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").Value = "65.86"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "236.47"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").Value = "2.50"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "23.42"
$ws.Range("E28").Value = "  +2.97%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.36"
$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("D31").Value = "156.78"
$ws.Range("E31").Value = "  +0.43%  "

$ws.Range("D32").Value = "32.46"
$ws.Range("E32").Value = "  +2.08%  "

$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("D35").Value = "3.02"
$ws.Range("E35").Value = "  +3.93%  "

$ws.Range("D36").Value = "0.0718"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("E38").Value = "  +1.82%  "

$ws.Range("E39").Value = "  +3.94%  "

$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  +3.13%  "

$ws.Range("D41").Value = "15.72"
$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("D42").Value = "3.85"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").Value = "2.067.28"
$ws.Range("E43").Value = "  -2.94%  "

$ws.Range("D44").Value = "19.46"
$ws.Range("E44").Value = "  +5.56%  "

$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").Value = "10.03"
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +7.15%  "

$ws.Range("E48").Value = "  -11.00%  "

$ws.Range("D49").Value = "2.430.43"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").Value = "1.48"
$ws.Range("E51").Value = "  +1.25%  "

